$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.67

# Row 3
$ws.Range("F3").Value = 1.09
$ws.Range("G3").Value = 26
$ws.Range("H3").Value = 1.3
$ws.Range("I3").Value = 1.47
$ws.Range("P3").Value = 2.58
$ws.Range("Q3").Value = 1.39

# Row 4
$ws.Range("F4").Value = 2.06
$ws.Range("G4").Value = 2.36
$ws.Range("I4").Value = 4.6
$ws.Range("Q4").Value = 1.95

# Row 5
$ws.Range("O5").Value = 1.2
$ws.Range("Q5").Value = 1.62
$ws.Range("R5").Value = 1.62
$ws.Range("S5").Value = 2.52
$ws.Range("T5").Value = 1.69
$ws.Range("U5").Value = 2.36
$ws.Range("Z5").Value = 55
$ws.Range("AD5").Value = 25
$ws.Range("AE5").Value = 75
$ws.Range("AK5").Value = 16
$ws.Range("AL5").Value = 29

# Row 6
$ws.Range("K6").Value = 3.85
$ws.Range("N6").Value = 4.5
$ws.Range("R6").Value = 1.47
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 1.68
$ws.Range("U6").Value = 2.38
$ws.Range("X6").Value = 19
$ws.Range("AA6").Value = 32
$ws.Range("AC6").Value = 8.4
$ws.Range("AF6").Value = 26
$ws.Range("AG6").Value = 14.5
$ws.Range("AJ6").Value = 60
$ws.Range("AL6").Value = 44
$ws.Range("AN6").Value = 30
$ws.Range("AO6").Value = 15

# Row 7
$ws.Range("F7").Value = 5.3
$ws.Range("H7").Value = 1.7
$ws.Range("I7").Value = 1.73
$ws.Range("K7").Value = 4.4
$ws.Range("N7").Value = 4.3
$ws.Range("P7").Value = 2.12
$ws.Range("R7").Value = 1.44
$ws.Range("S7").Value = 3.15
$ws.Range("T7").Value = 1.83
$ws.Range("U7").Value = 2.12
$ws.Range("X7").Value = 18
$ws.Range("Y7").Value = 9.199999999999999
$ws.Range("AB7").Value = 20
$ws.Range("AF7").Value = 44
$ws.Range("AG7").Value = 21
$ws.Range("AH7").Value = 20
$ws.Range("AI7").Value = 34
$ws.Range("AN7").Value = 85
$ws.Range("AO7").Value = 9.6

# Row 8
$ws.Range("M8").Value = 1.09
$ws.Range("P8").Value = 1.75
$ws.Range("Y8").Value = 12.5
$ws.Range("AA8").Value = 1000
$ws.Range("AD8").Value = 17
$ws.Range("AE8").Value = 55
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 1000
$ws.Range("AO8").Value = 80

# Row 9
$ws.Range("F9").Value = 2.58
$ws.Range("G9").Value = 2.82
$ws.Range("H9").Value = 2.82
$ws.Range("I9").Value = 3.05
$ws.Range("J9").Value = 3.35

# Row 10
$ws.Range("F10").Value = 1.58
$ws.Range("H10").Value = 6.4
$ws.Range("J10").Value = 3.95

# Row 11
$ws.Range("F11").Value = 2.9
$ws.Range("G11").Value = 3.2
$ws.Range("H11").Value = 2.54
$ws.Range("I11").Value = 2.76
$ws.Range("K11").Value = 3.6
$ws.Range("Q11").Value = 2.06

# Row 12
$ws.Range("F12").Value = 1.87
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 4.3
$ws.Range("J12").Value = 3.6

# Row 13
$ws.Range("F13").Value = 1.45
$ws.Range("H13").Value = 6.4
$ws.Range("I13").Value = 8

# Row 15
$ws.Range("F15").Value = 1.46
$ws.Range("G15").Value = 1.6
$ws.Range("H15").Value = 6.6
$ws.Range("K15").Value = 5.4

# Row 16
$ws.Range("G16").Value = 2.9
$ws.Range("I16").Value = 3.7

# Row 17
$ws.Range("H17").Value = 1.54
$ws.Range("I17").Value = 1.69
$ws.Range("J17").Value = 4.6

# Row 18
$ws.Range("Q18").Value = 1.72

# Row 20
$ws.Range("F20").Value = 2.48
$ws.Range("G20").Value = 3.15
$ws.Range("H20").Value = 2.84
$ws.Range("I20").Value = 3.65
$ws.Range("J20").Value = 2.82
$ws.Range("K20").Value = 3.65
$ws.Range("P20").Value = 1.73

# Row 21
$ws.Range("G21").Value = 1.69
$ws.Range("H21").Value = 6.2
$ws.Range("I21").Value = 8.6
$ws.Range("K21").Value = 5.5

# Row 22
$ws.Range("F22").Value = 2.22

# Row 23
$ws.Range("F23").Value = 1.86
$ws.Range("P23").Value = 1.81
$ws.Range("Q23").Value = 1.99

# Row 25
$ws.Range("G25").Value = 1.4
$ws.Range("P25").Value = 2.52
$ws.Range("T25").Value = 1.94
$ws.Range("AA25").Value = 390
$ws.Range("AB25").Value = 10
$ws.Range("AL25").Value = 34

# Row 26
$ws.Range("N26").Value = 3.9
$ws.Range("O26").Value = 1.32
$ws.Range("P26").Value = 2
$ws.Range("Q26").Value = 1.97
$ws.Range("R26").Value = 1.37
$ws.Range("S26").Value = 3.5
$ws.Range("T26").Value = 2.18
$ws.Range("AB26").Value = 7.8
$ws.Range("AF26").Value = 8
$ws.Range("AI26").Value = 160
$ws.Range("AJ26").Value = 12.5
